$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION1")

$ws.Range("D17").ClearContents()
$ws.Range("F17").Value = "x"

$ws.Range("D18").ClearContents()
$ws.Range("F18").Value = "x"
